$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Ovens"
$ws.Range("D5").Value = "https://www.samsung.com/in/microwave-ovens/all-microwave-ovens/"
$ws.Range("C5").Value = "All"
$ws.Range("B5").Value = "Samsung"

$ws.Columns.Item(4).ColumnWidth = 59.6

$ws.Range("G14").Select()
